$wb = $excel.ActiveWorkbook

# "BPCList" sheet: update the ExecFlag for RightsManagement (row 5) from Yes to No
$ws = $wb.Worksheets.Item("BPCList")
$ws.Activate()
$ws.Range("B5").Value = "No"

# reflect the new active selection on the sheet after the edit
$ws.Range("C10").Select()
